# Generate Report for Handback
#
# The localization-status report is regenerated once a handback package has
# come back in sync with en-US: the Status column is updated, the "Latest
# Target File" / "Latest Handback File" columns are (newly) populated with
# the file that was handed back, and "Latest Handback DateTime" moves from
# the zero-date placeholder to the real handback timestamp. Done for both
# the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$srcFileDisplay = "67741593-67a5-440a-95f7-1285f266c3b8.md"
$srcFileUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/99b87bf246d48f17aec2ec6d0ffd28f90bf0709e/e2e/67741593-67a5-440a-95f7-1285f266c3b8.md"

$configDisplay = ".localization-config"
$configUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/99b87bf246d48f17aec2ec6d0ffd28f90bf0709e/.localization-config"

function Set-ReportRow {
    param(
        $ws,
        $handoffDisplay,
        $handoffUrl,
        $handbackDateTime
    )

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = $newStatus

    # New columns: Latest Target File (E2) / Latest Handback File (F2)
    $ws.Range("E2").Value = $srcFileDisplay
    $ws.Range("F2").Value = $handoffDisplay

    # Latest Handback DateTime (G2): placeholder -> real timestamp
    $ws.Range("G2").Value = $handbackDateTime

    # The "Latest Handoff Datetime" column (D) keeps its datetime display
    # format; re-assert it on D2/D3 since re-saving the workbook otherwise
    # drops the number format that was already applied to those cells.
    $ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    # Rebuild the hyperlinks top-to-bottom / left-to-right so the
    # relationship ids come out in document order (A2, C2, E2, F2, A3).
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $srcFileUrl, [Type]::Missing, [Type]::Missing, $srcFileDisplay) | Out-Null
    $ws.Range("A2").Font.Underline = $true
    $ws.Range("A2").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("C2"), $handoffUrl, [Type]::Missing, [Type]::Missing, $handoffDisplay) | Out-Null
    $ws.Range("C2").Font.Underline = $true
    $ws.Range("C2").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("E2"), $srcFileUrl, [Type]::Missing, [Type]::Missing, $srcFileDisplay) | Out-Null
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("F2"), $handoffUrl, [Type]::Missing, [Type]::Missing, $handoffDisplay) | Out-Null
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, $configDisplay) | Out-Null
    $ws.Range("A3").Font.Underline = $true
    $ws.Range("A3").Font.Color = 15570276
}

# --- zh-cn sheet -------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhHandoffDisplay = "67741593-67a5-440a-95f7-1285f266c3b8.800dc078353e2829d40d817f958a4bd09af759e9.zh-cn.xlf"
$zhHandoffUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/05d24ce0f1f6e8239d7a7b2ea741331ed285587b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/67741593-67a5-440a-95f7-1285f266c3b8.800dc078353e2829d40d817f958a4bd09af759e9.zh-cn.xlf"

Set-ReportRow $wsZh $zhHandoffDisplay $zhHandoffUrl "2016-02-25 09:05:41"

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deHandoffDisplay = "67741593-67a5-440a-95f7-1285f266c3b8.800dc078353e2829d40d817f958a4bd09af759e9.de-de.xlf"
$deHandoffUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3c635007f7b8b685a9f0a6c4d0e3481b363f37f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/67741593-67a5-440a-95f7-1285f266c3b8.800dc078353e2829d40d817f958a4bd09af759e9.de-de.xlf"

Set-ReportRow $wsDe $deHandoffDisplay $deHandoffUrl "2016-02-25 09:06:00"
